$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.00449775112443778
$ws.Range("C2").Value = 0.00224887556221889
$ws.Range("D2").Value = 0.00524737631184408
$ws.Range("E2").Value = 0.00224887556221889
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.991754122938531
$ws.Range("H2").Value = 0.0217391304347826
$ws.Range("I2").Value = 0.995502248875562
$ws.Range("J2").Value = 0.0194902548725637
$ws.Range("K2").Value = 0.00449775112443778
$ws.Range("L2").Value = 0.0157421289355322
$ws.Range("M2").Value = 0.000749625187406297
$ws.Range("N2").Value = 0.027736131934033
$ws.Range("O2").Value = 0.00149925037481259
$ws.Range("P2").Value = 0.991754122938531
$ws.Range("Q2").Value = 0.00524737631184408
$ws.Range("R2").Value = 0.992503748125937
$ws.Range("S2").Value = 0.0434782608695652
$ws.Range("T2").Value = 0.991754122938531
$ws.Range("U2").Value = 0.991004497751124
$ws.Range("V2").Value = 0.00149925037481259
$ws.Range("W2").Value = 0.00299850074962519
$ws.Range("X2").Value = 0.986506746626687

$ws.Range("B3").Value = 0.986506746626687
$ws.Range("C3").Value = 0.992503748125937
$ws.Range("D3").Value = 0.991754122938531
$ws.Range("E3").Value = 0.00224887556221889
$ws.Range("F3").Value = 0.997751124437781
$ws.Range("G3").Value = 0.00149925037481259
$ws.Range("H3").Value = 0.00224887556221889
$ws.Range("I3").Value = 0.00374812593703148
$ws.Range("J3").Value = 0.00224887556221889
$ws.Range("K3").Value = 0.0217391304347826
$ws.Range("L3").Value = 0.979010494752624
$ws.Range("M3").Value = 0.00224887556221889
$ws.Range("N3").Value = 0.00299850074962519
$ws.Range("O3").Value = 0.991754122938531
$ws.Range("P3").Value = 0.00224887556221889
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.00149925037481259
$ws.Range("S3").Value = 0.00449775112443778
$ws.Range("T3").Value = 0.00149925037481259
$ws.Range("U3").Value = 0.00599700149925037
$ws.Range("V3").Value = 0.00449775112443778
$ws.Range("W3").Value = 0.989505247376312
$ws.Range("X3").Value = 0.00299850074962519

$ws.Range("B4").Value = 0.00374812593703148
$ws.Range("C4").Value = 0.00449775112443778
$ws.Range("D4").Value = 0.00224887556221889
$ws.Range("E4").Value = 0.00149925037481259
$ws.Range("F4").Value = 0.000749625187406297
$ws.Range("G4").Value = 0.00299850074962519
$ws.Range("H4").Value = 0.971514242878561
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0.00449775112443778
$ws.Range("K4").Value = 0.97376311844078
$ws.Range("L4").Value = 0.00224887556221889
$ws.Range("M4").Value = 0.00149925037481259
$ws.Range("N4").Value = 0.968515742128936
$ws.Range("O4").Value = 0.00449775112443778
$ws.Range("P4").Value = 0.00374812593703148
$ws.Range("Q4").Value = 0.00299850074962519
$ws.Range("R4").Value = 0.00449775112443778
$ws.Range("S4").Value = 0.00149925037481259
$ws.Range("T4").Value = 0.00224887556221889
$ws.Range("U4").Value = 0.00224887556221889
$ws.Range("V4").Value = 0.988755622188906
$ws.Range("W4").Value = 0.00374812593703148
$ws.Range("X4").Value = 0.00524737631184408

$ws.Range("B5").Value = 0.00524737631184408
$ws.Range("C5").Value = 0.000749625187406297
$ws.Range("D5").Value = 0.000749625187406297
$ws.Range("E5").Value = 0.99400299850075
$ws.Range("F5").Value = 0.00149925037481259
$ws.Range("G5").Value = 0.00374812593703148
$ws.Range("H5").Value = 0.00449775112443778
$ws.Range("I5").Value = 0.000749625187406297
$ws.Range("J5").Value = 0.97376311844078
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0.00224887556221889
$ws.Range("M5").Value = 0.995502248875562
$ws.Range("N5").Value = 0.000749625187406297
$ws.Range("O5").Value = 0.00224887556221889
$ws.Range("P5").Value = 0.00224887556221889
$ws.Range("Q5").Value = 0.991754122938531
$ws.Range("R5").Value = 0.00149925037481259
$ws.Range("S5").Value = 0.950524737631184
$ws.Range("T5").Value = 0.00449775112443778
$ws.Range("U5").Value = 0.000749625187406297
$ws.Range("V5").Value = 0.00524737631184408
$ws.Range("W5").Value = 0.00374812593703148
$ws.Range("X5").Value = 0.00524737631184408
